$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 500010  # H6: 333346.66 -> 500010
$ws.Cells.Item(6, 9).Value = 500010  # I6: 333346.66 -> 500010
$ws.Cells.Item(6, 11).Value = 1500030  # K6: 1000039.98 -> 1500030
$ws.Cells.Item(6, 13).Value = -1499918  # M6: -999927.98 -> -1499918
$ws.Cells.Item(9, 8).Value = 396.64285  # H9: 375.4 -> 396.64285
$ws.Cells.Item(9, 9).Value = 429.41666  # I9: 398.53845 -> 429.41666
$ws.Cells.Item(9, 10).Value = 200  # J9: 225 -> 200
$ws.Cells.Item(9, 11).Value = 429.41666  # K9: 398.53845 -> 429.41666
$ws.Cells.Item(9, 12).Value = 200  # L9: 225 -> 200
$ws.Cells.Item(9, 13).Value = -260.41666  # M9: -229.53845 -> -260.41666
$ws.Cells.Item(9, 14).Value = -538  # N9: -563 -> -538
$ws.Cells.Item(12, 8).Value = 0  # H12: 210.5 -> 0
$ws.Cells.Item(12, 9).Value = 0  # I12: 210.5 -> 0
$ws.Cells.Item(12, 11).Value = 0  # K12: 210.5 -> 0
$ws.Cells.Item(12, 13).ClearContents()  # M12: was -40.5
$ws.Cells.Item(70, 8).Value = 2341.9167  # H70: 2207.923 -> 2341.9167
$ws.Cells.Item(70, 9).Value = 1637.5  # I70: 1522.2222 -> 1637.5
$ws.Cells.Item(70, 11).Value = 4912.5  # K70: 4566.6666 -> 4912.5
$ws.Cells.Item(70, 13).Value = -4642.5  # M70: -4296.6666 -> -4642.5
$ws.Cells.Item(73, 8).Value = 2341.9167  # H73: 2207.923 -> 2341.9167
$ws.Cells.Item(73, 9).Value = 1637.5  # I73: 1522.2222 -> 1637.5
$ws.Cells.Item(73, 11).Value = 4912.5  # K73: 4566.6666 -> 4912.5
$ws.Cells.Item(73, 13).Value = -3976.5  # M73: -3630.6666 -> -3976.5
$ws.Cells.Item(116, 8).Value = 4498.5  # H116: 4465.6665 -> 4498.5
$ws.Cells.Item(116, 9).Value = 4498  # I116: 4449 -> 4498
$ws.Cells.Item(116, 11).Value = 4498  # K116: 4449 -> 4498
$ws.Cells.Item(116, 13).Value = -1056  # M116: -1007 -> -1056

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1905.9412  # H2: 1598.5385 -> 1905.9412
$ws.Cells.Item(2, 9).Value = 1530.8462  # I2: 1336.3334 -> 1530.8462
$ws.Cells.Item(2, 10).Value = 3125  # J2: 2699.8 -> 3125
$ws.Cells.Item(2, 11).Value = 1530.8462  # K2: 1336.3334 -> 1530.8462
$ws.Cells.Item(2, 12).Value = 3125  # L2: 2699.8 -> 3125
$ws.Cells.Item(2, 13).Value = -1417.8462  # M2: -1223.3334 -> -1417.8462
$ws.Cells.Item(2, 14).Value = -3351  # N2: -2925.8 -> -3351
$ws.Cells.Item(5, 8).Value = 491.42856  # H5: 492.57144 -> 491.42856
$ws.Cells.Item(5, 9).Value = 286  # I5: 286.1111 -> 286
$ws.Cells.Item(5, 10).Value = 861.2  # J5: 864.2 -> 861.2
$ws.Cells.Item(5, 11).Value = 286  # K5: 286.1111 -> 286
$ws.Cells.Item(5, 12).Value = 861.2  # L5: 864.2 -> 861.2
$ws.Cells.Item(5, 13).Value = -174  # M5: -174.1111 -> -174
$ws.Cells.Item(5, 14).Value = -1085.2  # N5: -1088.2 -> -1085.2
$ws.Cells.Item(29, 8).Value = 3599.5  # H29: 3633.3333 -> 3599.5
$ws.Cells.Item(29, 9).Value = 3599.5  # I29: 3633.3333 -> 3599.5
$ws.Cells.Item(29, 11).Value = 3599.5  # K29: 3633.3333 -> 3599.5
$ws.Cells.Item(29, 13).Value = -3291.5  # M29: -3325.3333 -> -3291.5
$ws.Cells.Item(32, 8).Value = 4088.027  # H32: 4194.0835 -> 4088.027
$ws.Cells.Item(32, 9).Value = 3036.0286  # I32: 3117.3823 -> 3036.0286
$ws.Cells.Item(32, 11).Value = 3036.0286  # K32: 3117.3823 -> 3036.0286
$ws.Cells.Item(32, 13).Value = -2749.0286  # M32: -2830.3823 -> -2749.0286
$ws.Cells.Item(60, 8).Value = 34000  # H60: 32200 -> 34000
$ws.Cells.Item(60, 10).Value = 0  # J60: 25000 -> 0
$ws.Cells.Item(60, 12).Value = 0  # L60: 25000 -> 0
$ws.Cells.Item(60, 14).ClearContents()  # N60: was -26466
$ws.Cells.Item(61, 8).Value = 1862.5  # H61: 1848.6666 -> 1862.5
$ws.Cells.Item(61, 9).Value = 1862.5  # I61: 1848.6666 -> 1862.5
$ws.Cells.Item(61, 11).Value = 1862.5  # K61: 1848.6666 -> 1862.5
$ws.Cells.Item(61, 13).Value = -1650.5  # M61: -1636.6666 -> -1650.5
$ws.Cells.Item(110, 8).Value = 499  # H110: 0 -> 499
$ws.Cells.Item(110, 9).Value = 499  # I110: 0 -> 499
$ws.Cells.Item(110, 11).Value = 499  # K110: 0 -> 499
$ws.Cells.Item(110, 13).Value = 1546  # M110: None -> 1546
$ws.Cells.Item(116, 8).Value = 1905.9412  # H116: 1598.5385 -> 1905.9412
$ws.Cells.Item(116, 9).Value = 1530.8462  # I116: 1336.3334 -> 1530.8462
$ws.Cells.Item(116, 10).Value = 3125  # J116: 2699.8 -> 3125
$ws.Cells.Item(116, 11).Value = 1530.8462  # K116: 1336.3334 -> 1530.8462
$ws.Cells.Item(116, 12).Value = 3125  # L116: 2699.8 -> 3125
$ws.Cells.Item(116, 13).Value = 763.1538  # M116: 957.6666 -> 763.1538
$ws.Cells.Item(116, 14).Value = -7713  # N116: -7287.8 -> -7713
$ws.Cells.Item(136, 8).Value = 1862.5  # H136: 1848.6666 -> 1862.5
$ws.Cells.Item(136, 9).Value = 1862.5  # I136: 1848.6666 -> 1862.5
$ws.Cells.Item(136, 11).Value = 5587.5  # K136: 5545.9998 -> 5587.5
$ws.Cells.Item(136, 13).Value = -3037.5  # M136: -2995.9998 -> -3037.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1905.9412  # H3: 1598.5385 -> 1905.9412
$ws.Cells.Item(3, 9).Value = 1530.8462  # I3: 1336.3334 -> 1530.8462
$ws.Cells.Item(3, 10).Value = 3125  # J3: 2699.8 -> 3125
$ws.Cells.Item(3, 11).Value = 1530.8462  # K3: 1336.3334 -> 1530.8462
$ws.Cells.Item(3, 12).Value = 3125  # L3: 2699.8 -> 3125
$ws.Cells.Item(3, 13).Value = -1416.8462  # M3: -1222.3334 -> -1416.8462
$ws.Cells.Item(3, 14).Value = -3353  # N3: -2927.8 -> -3353
$ws.Cells.Item(4, 8).Value = 491.42856  # H4: 492.57144 -> 491.42856
$ws.Cells.Item(4, 9).Value = 286  # I4: 286.1111 -> 286
$ws.Cells.Item(4, 10).Value = 861.2  # J4: 864.2 -> 861.2
$ws.Cells.Item(4, 11).Value = 286  # K4: 286.1111 -> 286
$ws.Cells.Item(4, 12).Value = 861.2  # L4: 864.2 -> 861.2
$ws.Cells.Item(4, 13).Value = -171  # M4: -171.1111 -> -171
$ws.Cells.Item(4, 14).Value = -1091.2  # N4: -1094.2 -> -1091.2
$ws.Cells.Item(86, 8).Value = 13904.454  # H86: 13904.909 -> 13904.454
$ws.Cells.Item(86, 9).Value = 14572.167  # I86: 14572.723 -> 14572.167
$ws.Cells.Item(86, 11).Value = 14572.167  # K86: 14572.723 -> 14572.167
$ws.Cells.Item(86, 13).Value = -13449.167  # M86: -13449.723 -> -13449.167
$ws.Cells.Item(89, 8).Value = 13904.454  # H89: 13904.909 -> 13904.454
$ws.Cells.Item(89, 9).Value = 14572.167  # I89: 14572.723 -> 14572.167
$ws.Cells.Item(89, 11).Value = 72860.83499999999  # K89: 72863.61500000001 -> 72860.83499999999
$ws.Cells.Item(89, 13).Value = -67244.83499999999  # M89: -67247.61500000001 -> -67244.83499999999
$ws.Cells.Item(102, 8).Value = 0  # H102: 20500 -> 0
$ws.Cells.Item(102, 9).Value = 0  # I102: 20500 -> 0
$ws.Cells.Item(102, 11).Value = 0  # K102: 20500 -> 0
$ws.Cells.Item(102, 13).ClearContents()  # M102: was -17255
$ws.Cells.Item(111, 8).Value = 10000  # H111: 24500 -> 10000
$ws.Cells.Item(111, 10).Value = 10000  # J111: 24500 -> 10000
$ws.Cells.Item(111, 12).Value = 10000  # L111: 24500 -> 10000
$ws.Cells.Item(111, 14).Value = -18180  # N111: -32680 -> -18180
$ws.Cells.Item(134, 8).Value = 1754.8948  # H134: 1782.75 -> 1754.8948
$ws.Cells.Item(134, 9).Value = 1740.5  # I134: 1782.75 -> 1740.5
$ws.Cells.Item(134, 10).Value = 2014  # J134: 0 -> 2014
$ws.Cells.Item(134, 11).Value = 5221.5  # K134: 5348.25 -> 5221.5
$ws.Cells.Item(134, 12).Value = 6042  # L134: 0 -> 6042
$ws.Cells.Item(134, 13).Value = -2686.5  # M134: -2813.25 -> -2686.5
$ws.Cells.Item(134, 14).Value = -11112  # N134: None -> -11112

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 34482930  # H7: 35714452 -> 34482930
$ws.Cells.Item(7, 9).Value = 38461680  # I7: 37037172 -> 38461680
$ws.Cells.Item(7, 10).Value = 416.33334  # J7: 999 -> 416.33334
$ws.Cells.Item(7, 11).Value = 38461680  # K7: 37037172 -> 38461680
$ws.Cells.Item(7, 12).Value = 416.33334  # L7: 999 -> 416.33334
$ws.Cells.Item(7, 13).Value = -38461567  # M7: -37037059 -> -38461567
$ws.Cells.Item(7, 14).Value = -642.33334  # N7: -1225 -> -642.33334
$ws.Cells.Item(60, 8).Value = 28012.75  # H60: 23450.375 -> 28012.75
$ws.Cells.Item(60, 9).Value = 4333.3335  # I60: 3625 -> 4333.3335
$ws.Cells.Item(60, 10).Value = 42220.4  # J60: 43275.75 -> 42220.4
$ws.Cells.Item(60, 11).Value = 4333.3335  # K60: 3625 -> 4333.3335
$ws.Cells.Item(60, 12).Value = 42220.4  # L60: 43275.75 -> 42220.4
$ws.Cells.Item(60, 13).Value = -3822.3335  # M60: -3114 -> -3822.3335
$ws.Cells.Item(60, 14).Value = -43242.4  # N60: -44297.75 -> -43242.4
$ws.Cells.Item(86, 8).Value = 11246.2  # H86: 11141.5 -> 11246.2
$ws.Cells.Item(86, 9).Value = 12018  # I86: 12118.857 -> 12018
$ws.Cells.Item(86, 11).Value = 12018  # K86: 12118.857 -> 12018
$ws.Cells.Item(86, 13).Value = -10895  # M86: -10995.857 -> -10895
$ws.Cells.Item(89, 8).Value = 11246.2  # H89: 11141.5 -> 11246.2
$ws.Cells.Item(89, 9).Value = 12018  # I89: 12118.857 -> 12018
$ws.Cells.Item(89, 11).Value = 60090  # K89: 60594.285 -> 60090
$ws.Cells.Item(89, 13).Value = -54474  # M89: -54978.285 -> -54474
$ws.Cells.Item(105, 8).Value = 0  # H105: 2241.5 -> 0
$ws.Cells.Item(105, 9).Value = 0  # I105: 2241.5 -> 0
$ws.Cells.Item(105, 11).Value = 0  # K105: 2241.5 -> 0
$ws.Cells.Item(105, 13).ClearContents()  # M105: was -494.5
$ws.Cells.Item(107, 8).Value = 1211.7778  # H107: 1200.6 -> 1211.7778
$ws.Cells.Item(107, 9).Value = 1211.7778  # I107: 1158.5264 -> 1211.7778
$ws.Cells.Item(107, 10).Value = 0  # J107: 2000 -> 0
$ws.Cells.Item(107, 11).Value = 1211.7778  # K107: 1158.5264 -> 1211.7778
$ws.Cells.Item(107, 12).Value = 0  # L107: 2000 -> 0
$ws.Cells.Item(107, 13).Value = 708.2221999999999  # M107: 761.4736 -> 708.2221999999999
$ws.Cells.Item(107, 14).ClearContents()  # N107: was -5840
$ws.Cells.Item(134, 8).Value = 2150  # H134: 1451.6 -> 2150
$ws.Cells.Item(134, 9).Value = 2150  # I134: 1749 -> 2150
$ws.Cells.Item(134, 10).Value = 0  # J134: 1005.5 -> 0
$ws.Cells.Item(134, 11).Value = 6450  # K134: 5247 -> 6450
$ws.Cells.Item(134, 12).Value = 0  # L134: 3016.5 -> 0
$ws.Cells.Item(134, 13).Value = -3915  # M134: -2712 -> -3915
$ws.Cells.Item(134, 14).ClearContents()  # N134: was -8086.5
$ws.Cells.Item(141, 8).Value = 329480.88  # H141: 303425.78 -> 329480.88
$ws.Cells.Item(141, 10).Value = 329480.88  # J141: 303425.78 -> 329480.88
$ws.Cells.Item(141, 12).Value = 329480.88  # L141: 303425.78 -> 329480.88
$ws.Cells.Item(141, 14).Value = -339840.88  # N141: -313785.78 -> -339840.88

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 1463  # H2: 2525 -> 1463
$ws.Cells.Item(2, 10).Value = 1934  # J2: 5000 -> 1934
$ws.Cells.Item(2, 12).Value = 11604  # L2: 30000 -> 11604
$ws.Cells.Item(2, 14).Value = -11830  # N2: -30226 -> -11830
$ws.Cells.Item(17, 8).Value = 783.3333  # H17: 486 -> 783.3333
$ws.Cells.Item(17, 9).Value = 100  # I17: 60 -> 100
$ws.Cells.Item(17, 11).Value = 300  # K17: 180 -> 300
$ws.Cells.Item(17, 13).Value = -131  # M17: -11 -> -131
$ws.Cells.Item(56, 8).Value = 9470.799999999999  # H56: 9790.143 -> 9470.799999999999
$ws.Cells.Item(56, 9).Value = 9470.799999999999  # I56: 9790.143 -> 9470.799999999999
$ws.Cells.Item(56, 11).Value = 9470.799999999999  # K56: 9790.143 -> 9470.799999999999
$ws.Cells.Item(56, 13).Value = -8940.799999999999  # M56: -9260.143 -> -8940.799999999999
$ws.Cells.Item(108, 8).Value = 885.7778  # H108: 979.125 -> 885.7778
$ws.Cells.Item(108, 9).Value = 1098.8572  # I108: 1258.8334 -> 1098.8572
$ws.Cells.Item(108, 11).Value = 3296.5716  # K108: 3776.5002 -> 3296.5716
$ws.Cells.Item(108, 13).Value = -416.5715999999998  # M108: -896.5001999999999 -> -416.5715999999998
$ws.Cells.Item(140, 8).Value = 591257.25  # H140: 717788.0600000001 -> 591257.25
$ws.Cells.Item(140, 9).Value = 591257.25  # I140: 717788.0600000001 -> 591257.25
$ws.Cells.Item(140, 11).Value = 1773771.75  # K140: 2153364.18 -> 1773771.75
$ws.Cells.Item(140, 13).Value = -1768591.75  # M140: -2148184.18 -> -1768591.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 8965.666999999999  # H55: 9099.75 -> 8965.666999999999
$ws.Cells.Item(55, 9).Value = 8965.666999999999  # I55: 9099.75 -> 8965.666999999999
$ws.Cells.Item(55, 11).Value = 8965.666999999999  # K55: 9099.75 -> 8965.666999999999
$ws.Cells.Item(55, 13).Value = -8638.666999999999  # M55: -8772.75 -> -8638.666999999999
$ws.Cells.Item(113, 8).Value = 2474.4  # H113: 2639.5557 -> 2474.4
$ws.Cells.Item(113, 9).Value = 2474.4  # I113: 2639.5557 -> 2474.4
$ws.Cells.Item(113, 11).Value = 2474.4  # K113: 2639.5557 -> 2474.4
$ws.Cells.Item(113, 13).Value = -304.4000000000001  # M113: -469.5556999999999 -> -304.4000000000001
$ws.Cells.Item(128, 8).Value = 106979.2  # H128: 106980 -> 106979.2
$ws.Cells.Item(128, 10).Value = 106979.2  # J128: 106980 -> 106979.2
$ws.Cells.Item(128, 12).Value = 106979.2  # L128: 106980 -> 106979.2
$ws.Cells.Item(128, 14).Value = -116939.2  # N128: -116940 -> -116939.2
$ws.Cells.Item(132, 8).Value = 15875206  # H132: 15875216 -> 15875206
$ws.Cells.Item(132, 9).Value = 1624.4375  # I132: 1637.5625 -> 1624.4375
$ws.Cells.Item(132, 11).Value = 4873.3125  # K132: 4912.6875 -> 4873.3125
$ws.Cells.Item(132, 13).Value = -2343.3125  # M132: -2382.6875 -> -2343.3125
$ws.Cells.Item(141, 8).Value = 129052.5  # H141: 129219.164 -> 129052.5
$ws.Cells.Item(141, 10).Value = 129052.5  # J141: 129219.164 -> 129052.5
$ws.Cells.Item(141, 12).Value = 129052.5  # L141: 129219.164 -> 129052.5
$ws.Cells.Item(141, 14).Value = -139412.5  # N141: -139579.164 -> -139412.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 35000  # H42: 0 -> 35000
$ws.Cells.Item(42, 10).Value = 35000  # J42: 0 -> 35000
$ws.Cells.Item(42, 12).Value = 35000  # L42: 0 -> 35000
$ws.Cells.Item(42, 14).Value = -36126  # N42: None -> -36126
$ws.Cells.Item(46, 8).Value = 3818.7646  # H46: 4010.5625 -> 3818.7646
$ws.Cells.Item(46, 9).Value = 1581.2858  # I46: 1719.8334 -> 1581.2858
$ws.Cells.Item(46, 11).Value = 1581.2858  # K46: 1719.8334 -> 1581.2858
$ws.Cells.Item(46, 13).Value = -1393.2858  # M46: -1531.8334 -> -1393.2858
$ws.Cells.Item(49, 8).Value = 35000  # H49: 0 -> 35000
$ws.Cells.Item(49, 10).Value = 35000  # J49: 0 -> 35000
$ws.Cells.Item(49, 12).Value = 35000  # L49: 0 -> 35000
$ws.Cells.Item(49, 14).Value = -35294  # N49: None -> -35294
$ws.Cells.Item(106, 8).Value = 16789.4  # H106: 17038 -> 16789.4
$ws.Cells.Item(106, 10).Value = 16789.4  # J106: 17038 -> 16789.4
$ws.Cells.Item(106, 12).Value = 16789.4  # L106: 17038 -> 16789.4
$ws.Cells.Item(106, 14).Value = -19313.4  # N106: -19562 -> -19313.4
$ws.Cells.Item(128, 8).Value = 73153.57000000001  # H128: 70255.5 -> 73153.57000000001
$ws.Cells.Item(128, 10).Value = 73153.57000000001  # J128: 70255.5 -> 73153.57000000001
$ws.Cells.Item(128, 12).Value = 73153.57000000001  # L128: 70255.5 -> 73153.57000000001
$ws.Cells.Item(128, 14).Value = -83113.57000000001  # N128: -80215.5 -> -83113.57000000001
$ws.Cells.Item(132, 8).Value = 1059  # H132: 1150 -> 1059
$ws.Cells.Item(132, 9).Value = 1074.25  # I132: 1150 -> 1074.25
$ws.Cells.Item(132, 10).Value = 998  # J132: 0 -> 998
$ws.Cells.Item(132, 11).Value = 3222.75  # K132: 3450 -> 3222.75
$ws.Cells.Item(132, 12).Value = 2994  # L132: 0 -> 2994
$ws.Cells.Item(132, 13).Value = -692.75  # M132: -920 -> -692.75
$ws.Cells.Item(132, 14).Value = -8054  # N132: None -> -8054

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 31333  # H103: 49000 -> 31333
$ws.Cells.Item(103, 10).Value = 31333  # J103: 49000 -> 31333
$ws.Cells.Item(103, 12).Value = 31333  # L103: 49000 -> 31333
$ws.Cells.Item(103, 14).Value = -33677  # N103: -51344 -> -33677
$ws.Cells.Item(124, 8).Value = 26616.75  # H124: 26618.25 -> 26616.75
$ws.Cells.Item(124, 10).Value = 26616.75  # J124: 26618.25 -> 26616.75
$ws.Cells.Item(124, 12).Value = 26616.75  # L124: 26618.25 -> 26616.75
$ws.Cells.Item(124, 14).Value = -36436.75  # N124: -36438.25 -> -36436.75
$ws.Cells.Item(135, 8).Value = 38567  # H135: 38425.25 -> 38567
$ws.Cells.Item(135, 10).Value = 38567  # J135: 38425.25 -> 38567
$ws.Cells.Item(135, 12).Value = 38567  # L135: 38425.25 -> 38567
$ws.Cells.Item(135, 14).Value = -48707  # N135: -48565.25 -> -48707
$ws.Cells.Item(140, 8).Value = 54759.43  # H140: 54902 -> 54759.43
$ws.Cells.Item(140, 10).Value = 54759.43  # J140: 54902 -> 54759.43
$ws.Cells.Item(140, 12).Value = 54759.43  # L140: 54902 -> 54759.43
$ws.Cells.Item(140, 14).Value = -65119.43  # N140: -65262 -> -65119.43
